$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update relative influence values for rows 2-5 (labels unchanged)
$ws.Range("B2").Value = 57.77627387403908
$ws.Range("B3").Value = 10.50512439712741
$ws.Range("B4").Value = 9.861220867180979
$ws.Range("B5").Value = 8.631631206992598

# Rows 6 and 7 swap variable order (Salinity now before SuspendedParticulateMatter)
# with updated relative influence values
$ws.Range("A6").Value = "Salinity"
$ws.Range("B6").Value = 6.63658334740203
$ws.Range("A7").Value = "SuspendedParticulateMatter"
$ws.Range("B7").Value = 6.589166307257903
